# ---------------------------------------------------------------------------
# Updates the dependencies deck:
#   1. Refresh the cached "today" date shown by the datetimeFigureOut field
#      on the slide master and every slide layout (9/25/2024 -> 9/26/2024).
#   2. On slide 2 (the ecommerce platform dependency diagram), shift the
#      "create_order" code label up to make room, and add a new
#      "get_order" code label below it (matching the existing style).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# EMU-per-point constant used by the PowerPoint COM object model (Shape.Left/
# Top/Width/Height and AddTextbox's placement args are all expressed in
# points). Adding a half-EMU fudge below compensates for the single-precision
# float the host stores these coordinates in, so the EMU value that ends up
# serialized to the XML lands exactly on the target integer instead of
# rounding down.
$emuPerPt = 12700.0
$halfEmu = 0.5 / $emuPerPt

function ToPt([double]$emu) {
    return ($emu / $emuPerPt) + $halfEmu
}

# --- 1. Re-stamp the cached date field wherever it appears -----------------

$newDate = "9/26/2024"

$master = $p.SlideMaster
for ($k = 1; $k -le $master.Shapes.Count; $k++) {
    $shp = $master.Shapes.Item($k)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2. Slide 2: reposition "create_order" and add "get_order" -------------

$slide2 = $p.Slides.Item(2)

$createOrder = $null
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shp = $slide2.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "create_order") {
        $createOrder = $shp
    }
}

$createOrder.Left = ToPt(3942232)
$createOrder.Top = ToPt(1335301)

$getOrder = $createOrder.Duplicate()
$getOrder.Name = "TextBox 1"
$getOrder.Left = ToPt(3949521)
$getOrder.Top = ToPt(1610638)
$getOrder.TextFrame.TextRange.Text = "get_order"
